$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows into the sprint board ---
# Current layout (rows 1-13):
#  1 header
#  2-4  SPRINT1
#  5-7  SPRINT2
#  8-9  SPRINT3 (view monthly bill / send bills)
#  10   SPRINT3 (copy over existing expenses)            -> needs a new SPRINT3 story before it
#  11-13 SPRINT4 (validate input / get avatar / dispay avatar) -> needs a new SPRINT4 story before "validate input"

# Insert a blank row at 10: pushes old row10 (copy expenses) down to row11,
# and old rows 11-13 (SPRINT4 stories) down to rows 12-14.
$ws.Rows.Item(10).Insert()

# After the previous insert, "validate input" (SPRINT4) now sits at row 12.
# Insert a blank row before it so the new SPRINT4 story can go at row 12.
$ws.Rows.Item(12).Insert()

# --- Fill the two newly inserted rows, copying formatting from a same-status row ---

# Row 10: new SPRINT3 / DONE story -> copy format+value from row 9 (SPRINT3/DONE), then overwrite the story text
$ws.Range("A9:C9").Copy($ws.Range("A10:C10"))
$ws.Range("B10").Value = "I want to view fee and bill for previous months"

# Row 12: new SPRINT4 / NOT STARTED story -> copy format+value from row 13 (SPRINT4/NOT STARTED, formerly row 11), then overwrite the story text
$ws.Range("A13:C13").Copy($ws.Range("A12:C12"))
$ws.Range("B12").Value = "I want to export and reimport DB"

# --- Resize the table to include the two new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C15"))

# --- Update the active selection to reflect where the edit was made ---
$ws.Range("C10").Select()
